$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("July 06, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "July 09, 2022", 2)
$find.Execute("September 04, 2022", $true, $false, $false, $false, $false, $true, 1, $false, "September 07, 2022", 2)
